$d = $word.ActiveDocument
$count = 0
$null = $d.Content.Find.Execute("40-27=13", $false, $false, $false, $false, $false, $true, 1, $false, "50+0=50", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("36+37=73", $false, $false, $false, $false, $false, $true, 1, $false, "71-17=54", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("94-27=67", $false, $false, $false, $false, $false, $true, 1, $false, "66+6=72", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("25-18=7", $false, $false, $false, $false, $false, $true, 1, $false, "89-71=18", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("44-23=21", $false, $false, $false, $false, $false, $true, 1, $false, "15+70=85", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("8+9=17", $false, $false, $false, $false, $false, $true, 1, $false, "4+58=62", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("32+48=80", $false, $false, $false, $false, $false, $true, 1, $false, "30-27=3", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("12+20=32", $false, $false, $false, $false, $false, $true, 1, $false, "32-27=5", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("56+34=90", $false, $false, $false, $false, $false, $true, 1, $false, "64+13=77", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("79-43=36", $false, $false, $false, $false, $false, $true, 1, $false, "68+23=91", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("61-60=1", $false, $false, $false, $false, $false, $true, 1, $false, "29+69=98", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("85-9=76", $false, $false, $false, $false, $false, $true, 1, $false, "6+77=83", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("15+26=41", $false, $false, $false, $false, $false, $true, 1, $false, "32-8=24", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("81-21=60", $false, $false, $false, $false, $false, $true, 1, $false, "58+3=61", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("72-30=42", $false, $false, $false, $false, $false, $true, 1, $false, "69-54=15", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("55-32=23", $false, $false, $false, $false, $false, $true, 1, $false, "88-31=57", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("80+13=93", $false, $false, $false, $false, $false, $true, 1, $false, "84-75=9", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("91-66=25", $false, $false, $false, $false, $false, $true, 1, $false, "82-54=28", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("8+40=48", $false, $false, $false, $false, $false, $true, 1, $false, "38+19=57", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("48+13=61", $false, $false, $false, $false, $false, $true, 1, $false, "50+18=68", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("48+10=58", $false, $false, $false, $false, $false, $true, 1, $false, "0+36=36", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("20+30=50", $false, $false, $false, $false, $false, $true, 1, $false, "11+76=87", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("82-23=59", $false, $false, $false, $false, $false, $true, 1, $false, "71-48=23", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("89-89=0", $false, $false, $false, $false, $false, $true, 1, $false, "40+47=87", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("45-1=44", $false, $false, $false, $false, $false, $true, 1, $false, "51-1=50", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("22+12=34", $false, $false, $false, $false, $false, $true, 1, $false, "83-52=31", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("2+62=64", $false, $false, $false, $false, $false, $true, 1, $false, "19+10=29", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("15+36=51", $false, $false, $false, $false, $false, $true, 1, $false, "73-0=73", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("59+26=85", $false, $false, $false, $false, $false, $true, 1, $false, "78+4=82", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("26+65=91", $false, $false, $false, $false, $false, $true, 1, $false, "29+14=43", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("56-11=45", $false, $false, $false, $false, $false, $true, 1, $false, "57-32=25", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("26-7=19", $false, $false, $false, $false, $false, $true, 1, $false, "15+68=83", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("12+84=96", $false, $false, $false, $false, $false, $true, 1, $false, "43+40=83", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("28+52=80", $false, $false, $false, $false, $false, $true, 1, $false, "64+30=94", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("7+68=75", $false, $false, $false, $false, $false, $true, 1, $false, "59+38=97", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("47+29=76", $false, $false, $false, $false, $false, $true, 1, $false, "28-15=13", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("69-28=41", $false, $false, $false, $false, $false, $true, 1, $false, "92-91=1", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("46+33=79", $false, $false, $false, $false, $false, $true, 1, $false, "13+67=80", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("78-59=19", $false, $false, $false, $false, $false, $true, 1, $false, "97-15=82", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("41-6=35", $false, $false, $false, $false, $false, $true, 1, $false, "76+5=81", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("95-66=29", $false, $false, $false, $false, $false, $true, 1, $false, "26-2=24", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("7+16=23", $false, $false, $false, $false, $false, $true, 1, $false, "35+60=95", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("98-34=64", $false, $false, $false, $false, $false, $true, 1, $false, "46+22=68", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("84-12=72", $false, $false, $false, $false, $false, $true, 1, $false, "54-24=30", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("9+1=10", $false, $false, $false, $false, $false, $true, 1, $false, "9+71=80", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("9+35=44", $false, $false, $false, $false, $false, $true, 1, $false, "18+8=26", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("30+0=30", $false, $false, $false, $false, $false, $true, 1, $false, "77-55=22", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("78-22=56", $false, $false, $false, $false, $false, $true, 1, $false, "72-4=68", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("8+70=78", $false, $false, $false, $false, $false, $true, 1, $false, "2+41=43", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("13+66=79", $false, $false, $false, $false, $false, $true, 1, $false, "64+35=99", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("8+64=72", $false, $false, $false, $false, $false, $true, 1, $false, "95-48=47", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("51-17=34", $false, $false, $false, $false, $false, $true, 1, $false, "32+18=50", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("44+22=66", $false, $false, $false, $false, $false, $true, 1, $false, "30+65=95", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("44-26=18", $false, $false, $false, $false, $false, $true, 1, $false, "58+32=90", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("60+12=72", $false, $false, $false, $false, $false, $true, 1, $false, "80+10=90", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("56+11=67", $false, $false, $false, $false, $false, $true, 1, $false, "94-71=23", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("84-83=1", $false, $false, $false, $false, $false, $true, 1, $false, "86-65=21", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("35+28=63", $false, $false, $false, $false, $false, $true, 1, $false, "46-14=32", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("19+79=98", $false, $false, $false, $false, $false, $true, 1, $false, "40-24=16", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("12+18=30", $false, $false, $false, $false, $false, $true, 1, $false, "53-2=51", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("74-42=32", $false, $false, $false, $false, $false, $true, 1, $false, "33-11=22", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("34+36=70", $false, $false, $false, $false, $false, $true, 1, $false, "80-77=3", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("86-20=66", $false, $false, $false, $false, $false, $true, 1, $false, "7+28=35", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("8+89=97", $false, $false, $false, $false, $false, $true, 1, $false, "94-7=87", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("22+18=40", $false, $false, $false, $false, $false, $true, 1, $false, "90-1=89", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("11-9=2", $false, $false, $false, $false, $false, $true, 1, $false, "66+0=66", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("35-33=2", $false, $false, $false, $false, $false, $true, 1, $false, "82-4=78", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("31+7=38", $false, $false, $false, $false, $false, $true, 1, $false, "25+19=44", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("84-52=32", $false, $false, $false, $false, $false, $true, 1, $false, "22+33=55", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("27-3=24", $false, $false, $false, $false, $false, $true, 1, $false, "5+79=84", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("95-4=91", $false, $false, $false, $false, $false, $true, 1, $false, "32-19=13", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("13+63=76", $false, $false, $false, $false, $false, $true, 1, $false, "15+1=16", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("7+6=13", $false, $false, $false, $false, $false, $true, 1, $false, "94-24=70", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("22-11=11", $false, $false, $false, $false, $false, $true, 1, $false, "7+37=44", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("59+29=88", $false, $false, $false, $false, $false, $true, 1, $false, "28+67=95", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("78-37=41", $false, $false, $false, $false, $false, $true, 1, $false, "59+38=97", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("24+35=59", $false, $false, $false, $false, $false, $true, 1, $false, "25+5=30", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("69-7=62", $false, $false, $false, $false, $false, $true, 1, $false, "48-4=44", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("48+4=52", $false, $false, $false, $false, $false, $true, 1, $false, "49-45=4", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("47+18=65", $false, $false, $false, $false, $false, $true, 1, $false, "24+29=53", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("51+16=67", $false, $false, $false, $false, $false, $true, 1, $false, "94-32=62", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("27+40=67", $false, $false, $false, $false, $false, $true, 1, $false, "92-22=70", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("64+32=96", $false, $false, $false, $false, $false, $true, 1, $false, "91-60=31", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("92-33=59", $false, $false, $false, $false, $false, $true, 1, $false, "3+79=82", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("13+86=99", $false, $false, $false, $false, $false, $true, 1, $false, "54-0=54", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("62-49=13", $false, $false, $false, $false, $false, $true, 1, $false, "24+60=84", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("38+17=55", $false, $false, $false, $false, $false, $true, 1, $false, "92-51=41", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("70-19=51", $false, $false, $false, $false, $false, $true, 1, $false, "56+21=77", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("83-6=77", $false, $false, $false, $false, $false, $true, 1, $false, "91-79=12", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("93+1=94", $false, $false, $false, $false, $false, $true, 1, $false, "97-69=28", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("50-19=31", $false, $false, $false, $false, $false, $true, 1, $false, "78-69=9", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("21+0=21", $false, $false, $false, $false, $false, $true, 1, $false, "42+34=76", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("63-7=56", $false, $false, $false, $false, $false, $true, 1, $false, "6+90=96", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("19+68=87", $false, $false, $false, $false, $false, $true, 1, $false, "69-24=45", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("49-22=27", $false, $false, $false, $false, $false, $true, 1, $false, "15+64=79", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("39+38=77", $false, $false, $false, $false, $false, $true, 1, $false, "13-11=2", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("63-30=33", $false, $false, $false, $false, $false, $true, 1, $false, "29+39=68", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("87-58=29", $false, $false, $false, $false, $false, $true, 1, $false, "59-41=18", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("48+29=77", $false, $false, $false, $false, $false, $true, 1, $false, "27+72=99", 2)
$count = $count + 1
$null = $d.Content.Find.Execute("24-16=8", $false, $false, $false, $false, $false, $true, 1, $false, "95-75=20", 2)
$count = $count + 1
Write-Output "Replacements attempted: $count"
